# Fruta / hortaliza, semanal
# Insert a new weekly data row for "Feria Lagunitas de Puerto Montt - Pomelo"
# at row 152, pushing existing rows 152:187 down to 153:188.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 152 (shifts rows 152-187 down to 153-188)
$ws.Rows.Item(152).Insert()

# Populate the new row 152 with the latest weekly observation
$ws.Range("A152").Value = 4
$ws.Range("B152").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C152").Value = "Los Lagos"
$ws.Range("D152").Value = 44543
$ws.Range("E152").Value = 10
$ws.Range("F152").Value = "Fruta"
$ws.Range("G152").Value = 100102
$ws.Range("H152").Value = "Cítricos"
$ws.Range("I152").Value = 100102006
$ws.Range("J152").Value = "Pomelo"
$ws.Range("K152").Value = "Start Ruby"
$ws.Range("L152").Value = "Primera"
$ws.Range("M152").Value = 120
$ws.Range("N152").Value = 11000
$ws.Range("O152").Value = 12000
$ws.Range("P152").Value = 11500
$ws.Range("Q152").Value = "`$/caja 14 kilos empedrada"
$ws.Range("R152").Value = "Región de O'Higgins"
$ws.Range("S152").Value = 821
$ws.Range("T152").Value = 14
